$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 0.2808873333333333
$ws.Range("H2").Value = 0.842662
$ws.Range("I2").Value = 0.5595554696739399
$ws.Range("J2").Value = 0.5595554696739399
$ws.Range("M2").Value = 29.17403400000001
$ws.Range("N2").Value = 87.52210200000002
$ws.Range("O2").Value = 0.3835306213274714
$ws.Range("P2").Value = 0.3835306213274714
$ws.Range("Q2").Value = 8.194616612836001
$ws.Range("R2").Value = 73.75154951552402
$ws.Range("S2").Value = 0.2146066569512312
$ws.Range("T2").Value = 0.2146066569512312

# Row 3
$ws.Range("G3").Value = 0.2808873333333333
$ws.Range("H3").Value = 0.842662
$ws.Range("I3").Value = 0.5595554696739399
$ws.Range("J3").Value = 0.5595554696739399
$ws.Range("O3").Value = 0.5274816184042599
$ws.Range("P3").Value = 0.5274816184042599
$ws.Range("Q3").Value = 11.27031165902778
$ws.Range("R3").Value = 101.43280493125
$ws.Range("S3").Value = 0.2951552247305655
$ws.Range("T3").Value = 0.2951552247305655

# Row 4
$ws.Range("G4").Value = 0.2808873333333333
$ws.Range("H4").Value = 0.842662
$ws.Range("I4").Value = 0.5595554696739399
$ws.Range("J4").Value = 0.5595554696739399
$ws.Range("M4").Value = 6.769034333333334
$ws.Range("N4").Value = 20.307103
$ws.Range("O4").Value = 0.08898776026826866
$ws.Range("P4").Value = 0.08898776026826867
$ws.Range("Q4").Value = 1.901336003131778
$ws.Range("R4").Value = 17.112024028186
$ws.Range("S4").Value = 0.04979358799214303
$ws.Range("T4").Value = 0.04979358799214304

# Row 5
$ws.Range("E5").Value = 1
$ws.Range("F5").Value = 0.3333333333333333
$ws.Range("G5").Value = 0.2210956666666667
$ws.Range("H5").Value = 0.663287
$ws.Range("I5").Value = 0.4404445303260602
$ws.Range("J5").Value = 0.4404445303260602
$ws.Range("M5").Value = 29.17403400000001
$ws.Range("N5").Value = 87.52210200000002
$ws.Range("O5").Value = 0.3835306213274714
$ws.Range("P5").Value = 0.3835306213274714
$ws.Range("Q5").Value = 6.450252496586002
$ws.Range("R5").Value = 58.05227246927401
$ws.Range("S5").Value = 0.1689239643762402
$ws.Range("T5").Value = 0.1689239643762402

# Row 6
$ws.Range("E6").Value = 1
$ws.Range("F6").Value = 0.3333333333333333
$ws.Range("G6").Value = 0.2210956666666667
$ws.Range("H6").Value = 0.663287
$ws.Range("I6").Value = 0.4404445303260602
$ws.Range("J6").Value = 0.4404445303260602
$ws.Range("O6").Value = 0.5274816184042599
$ws.Range("P6").Value = 0.5274816184042599
$ws.Range("Q6").Value = 8.871233317013889
$ws.Range("R6").Value = 79.841099853125
$ws.Range("S6").Value = 0.2323263936736944
$ws.Range("T6").Value = 0.2323263936736943

# Row 7
$ws.Range("E7").Value = 1
$ws.Range("F7").Value = 0.3333333333333333
$ws.Range("G7").Value = 0.2210956666666667
$ws.Range("H7").Value = 0.663287
$ws.Range("I7").Value = 0.4404445303260602
$ws.Range("J7").Value = 0.4404445303260602
$ws.Range("M7").Value = 6.769034333333334
$ws.Range("N7").Value = 20.307103
$ws.Range("O7").Value = 0.08898776026826866
$ws.Range("P7").Value = 0.08898776026826867
$ws.Range("Q7").Value = 1.496604158617889
$ws.Range("R7").Value = 13.469437427561
$ws.Range("S7").Value = 0.03919417227612564
$ws.Range("T7").Value = 0.03919417227612564
